$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared strings need to be created in the same order as the source
# workbook (name, name, reference) so the sharedStrings table matches:
#   PS J1011-01, PS J0202-17, Ighina+24
$ws.Range("A34").Value = "PS J1011-01"
$ws.Range("A35").Value = "PS J0202-17"

# Row 33: extend the existing J2020-6215 row with the new columns
$ws.Range("C33").Value = 29.5
$ws.Range("D33").Value = -0.27
$ws.Range("G33").Value = "Ighina+24"

# Row 34: new source PS J1011-01
$ws.Range("B34").Value = 5.58
$ws.Range("B34").HorizontalAlignment = -4152
$ws.Range("C34").Value = 7.46
$ws.Range("D34").Value = -0.4
$ws.Range("E34").Value = "Ighina+24"

# Row 35: new source PS J0202-17
$ws.Range("B35").Value = 5.57
$ws.Range("B35").HorizontalAlignment = -4152
$ws.Range("C35").Value = 43.16
$ws.Range("D35").Value = -0.65
$ws.Range("E35").Value = "Ighina+24"

# Update the selected cell / view state like Excel would after editing
$ws.Range("B34").Select()
